$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.817.82"
$ws.Range("E2").Value = "  -3.75%  "
$ws.Range("D3").Value = "1.815.70"
$ws.Range("E3").Value = "  -3.06%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'276.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -8.24%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "'0.5100"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.08%  "
$ws.Range("D8").Value = "'0.3522"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.35%  "
$ws.Range("D9").Value = "'44.72"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.80%  "
$ws.Range("D10").Value = "'0.06670"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.16%  "
$ws.Range("D11").Value = "'20.03"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.06%  "
$ws.Range("D12").Value = "'0.8300"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.51%  "
$ws.Range("D13").Value = "'0.07864"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.64%  "
$ws.Range("D14").Value = "1.806.14"
$ws.Range("E14").Value = "  -3.65%  "
$ws.Range("D15").Value = "'5.076"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.51%  "
$ws.Range("E16").Value = "  -6.47%  "
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").Value = "'14.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.08%  "
$ws.Range("D19").Value = "'0.000008022"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.18%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "25.866.45"
$ws.Range("E21").Value = "  -3.79%  "
$ws.Range("D22").Value = "'4.724"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.16%  "
$ws.Range("D23").Value = "'10.02"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.25%  "
$ws.Range("D24").Value = "'6.074"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.89%  "
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").Value = "'2.198"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.76%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'140.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.14%  "
$ws.Range("E27").Value = "  -3.65%  "
$ws.Range("E28").Value = "  -5.32%  "
$ws.Range("D29").Value = "'109.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.06%  "
$ws.Range("D30").Value = "'4.350"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.92%  "
$ws.Range("D31").Value = "'4.233"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.96%  "
$ws.Range("D32").Value = "'0.08799"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.91%  "
$ws.Range("D33").Value = "'0.04895"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.66%  "
$ws.Range("D34").Value = "'0.7312"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.09%  "
$ws.Range("D35").Value = "'1.137"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.97%  "
$ws.Range("D36").Value = "'2.881"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.34%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "'1.000"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'3.146"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("D39").Value = "'2.383"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.42%  "
$ws.Range("D40").Value = "'0.5199"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -13.58%  "
$ws.Range("D41").Value = "'0.01849"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.33%  "
$ws.Range("D42").Value = "'0.9563"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.89%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'6.208"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.53%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'111.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.31%  "
$ws.Range("D45").Value = "'8.019"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.17%  "
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").Value = "'0.4554"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -11.42%  "
$ws.Range("D48").Value = "'0.1362"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.69%  "
$ws.Range("D49").Value = "'36.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.60%  "
$ws.Range("D50").Value = "'9.266"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.44%  "
$ws.Range("E51").Value = "  -8.01%  "
